# Apply updated TPM-derived NATMI values to rows 2-13 (Icam1-Spn LR pairs).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2; A="ECs"; B="Icam1"; C="Spn"; D="ECs"; Nums=@(3, 1, 30.39114566666667, 91.173437, 0.2485034818803364, 0.2485034818803363, 1, 0.3333333333333333, 0.01649366666666667, 0.049481, 0.005209935740510855, 0.005209935740510855, 0.5012614262441112, 4.511352836197, 0.001294687171889756, 0.001294687171889756) },
    @{ Row=3; A="ECs"; B="Icam1"; C="Spn"; D="MuSCs"; Nums=@(3, 1, 30.39114566666667, 91.173437, 0.2485034818803364, 0.2485034818803363, 2, 0.6666666666666666, 0.133398, 0.400194, 0.04213708340045676, 0.04213708340045676, 4.054118049642, 36.487062446778, 0.01047121194129563, 0.01047121194129563) },
    @{ Row=4; A="ECs"; B="Icam1"; C="Spn"; D="Resolving-Mac"; Nums=@(3, 1, 30.39114566666667, 91.173437, 0.2485034818803364, 0.2485034818803363, 3, 1, 3.015918333333333, 9.047754999999999, 0.9526529808590324, 0.9526529808590325, 91.6572133871039, 824.9149204839349, 0.236737582767151, 0.236737582767151) },
    @{ Row=5; A="FAPs"; B="Icam1"; C="Spn"; D="ECs"; Nums=@(3, 1, 45.91529066666667, 137.745872, 0.3754419042757282, 0.3754419042757282, 1, 0.3333333333333333, 0.01649366666666667, 0.049481, 0.005209935740510855, 0.005209935740510855, 0.7573114991591112, 6.815803492432001, 0.001956028195571571, 0.001956028195571571) },
    @{ Row=6; A="FAPs"; B="Icam1"; C="Spn"; D="MuSCs"; Nums=@(3, 1, 45.91529066666667, 137.745872, 0.3754419042757282, 0.3754419042757282, 2, 0.6666666666666666, 0.133398, 0.400194, 0.04213708340045676, 0.04213708340045676, 6.125007944352, 55.12507149916801, 0.01582002683249266, 0.01582002683249266) },
    @{ Row=7; A="FAPs"; B="Icam1"; C="Spn"; D="Resolving-Mac"; Nums=@(3, 1, 45.91529066666667, 137.745872, 0.3754419042757282, 0.3754419042757282, 3, 1, 3.015918333333333, 9.047754999999999, 0.9526529808590324, 0.9526529808590325, 138.4767669019289, 1246.29090211736, 0.3576658492476639, 0.357665849247664) },
    @{ Row=8; A="MuSCs"; B="Icam1"; C="Spn"; D="ECs"; Nums=@(3, 1, 2.332475, 6.997425, 0.01907227076123622, 0.01907227076123622, 1, 0.3333333333333333, 0.01649366666666667, 0.049481, 0.005209935740510855, 0.005209935740510855, 0.03847106515833334, 0.346239586425, 0.00009936530509166474, 0.00009936530509166474) },
    @{ Row=9; A="MuSCs"; B="Icam1"; C="Spn"; D="MuSCs"; Nums=@(3, 1, 2.332475, 6.997425, 0.01907227076123622, 0.01907227076123622, 2, 0.6666666666666666, 0.133398, 0.400194, 0.04213708340045676, 0.04213708340045676, 0.31114750005, 2.80032750045, 0.0008036498637023034, 0.0008036498637023036) },
    @{ Row=10; A="MuSCs"; B="Icam1"; C="Spn"; D="Resolving-Mac"; Nums=@(3, 1, 2.332475, 6.997425, 0.01907227076123622, 0.01907227076123622, 3, 1, 3.015918333333333, 9.047754999999999, 0.9526529808590324, 0.9526529808590325, 7.034554114541666, 63.31098703087499, 0.01816925559244225, 0.01816925559244225) },
    @{ Row=11; A="Resolving-Mac"; B="Icam1"; C="Spn"; D="ECs"; Nums=@(3, 1, 43.657748, 130.973244, 0.3569823430826993, 0.3569823430826993, 1, 0.3333333333333333, 0.01649366666666667, 0.049481, 0.005209935740510855, 0.005209935740510855, 0.7200763429293333, 6.480687086363999, 0.001859855067957863, 0.001859855067957863) },
    @{ Row=12; A="Resolving-Mac"; B="Icam1"; C="Spn"; D="MuSCs"; Nums=@(3, 1, 43.657748, 130.973244, 0.3569823430826993, 0.3569823430826993, 2, 0.6666666666666666, 0.133398, 0.400194, 0.04213708340045676, 0.04213708340045676, 5.823856267704, 52.414706409336, 0.01504219476296617, 0.01504219476296617) },
    @{ Row=13; A="Resolving-Mac"; B="Icam1"; C="Spn"; D="Resolving-Mac"; Nums=@(3, 1, 43.657748, 130.973244, 0.3569823430826993, 0.3569823430826993, 3, 1, 3.015918333333333, 9.047754999999999, 0.9526529808590324, 0.9526529808590325, 131.6682025852467, 1185.01382326722, 0.3400802932517752, 0.3400802932517753) }
)

$numCols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Range("A$r").Value = $item.A
    $ws.Range("B$r").Value = $item.B
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D
    for ($i = 0; $i -lt $numCols.Count; $i++) {
        $ws.Range($numCols[$i] + "$r").Value = $item.Nums[$i]
    }
}

